# This script reproduces the cryptos-list data refresh described in the commit
# message ('Updated cryptos list ... with GitHub Actions'). All target cells
# already hold text (not numeric) values in the workbook, so we force the
# NumberFormat to "@" (Text) before writing any value that Excel could
# otherwise auto-convert/reformat as a number (e.g. "70.90" -> 70.9,
# "0.520" -> 0.52, "0.0741" -> 0.0741 as a float instead of text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.657.47'
$ws.Range("E2").Value = '  +2.40%  '
# Row 3
$ws.Range("D3").Value = '2.530.68'
# Row 4
$ws.Range("E4").Value = '  -0.03%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.99'
$ws.Range("E5").Value = '  +2.05%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.44'
$ws.Range("E6").Value = '  +2.01%  '
# Row 7
$ws.Range("E7").Value = '  -0.02%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.520'
$ws.Range("E8").Value = '  +1.69%  '
# Row 9
$ws.Range("D9").Value = '2.530.58'
$ws.Range("E9").Value = '  +2.65%  '
# Row 10
$ws.Range("E10").Value = '  +6.42%  '
# Row 11
$ws.Range("E11").Value = '  -1.00%  '
# Row 12
$ws.Range("E12").Value = '  +1.24%  '
# Row 13
$ws.Range("E13").Value = '  +1.97%  '
# Row 14
$ws.Range("D14").Value = '2.993.48'
$ws.Range("E14").Value = '  +2.76%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.19'
$ws.Range("E15").Value = '  +3.34%  '
# Row 16
$ws.Range("D16").Value = '68.603.07'
$ws.Range("E16").Value = '  +2.48%  '
# Row 17
$ws.Range("E17").Value = '  +1.21%  '
# Row 18
$ws.Range("D18").Value = '2.524.79'
$ws.Range("E18").Value = '  +2.86%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.09'
$ws.Range("E19").Value = '  +1.86%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.55'
$ws.Range("E20").Value = '  +1.47%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '352.97'
$ws.Range("E21").Value = '  +1.44%  '
# Row 22
$ws.Range("E22").Value = '  +5.11%  '
# Row 23
$ws.Range("E23").Value = '  +0.03%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.90'
$ws.Range("E24").Value = '  +2.35%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.25'
$ws.Range("E25").Value = '  +1.58%  '
# Row 26
$ws.Range("E26").Value = '  -4.64%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.04'
$ws.Range("E27").Value = '  -1.40%  '
# Row 28
$ws.Range("D28").Value = '2.690.83'
$ws.Range("E28").Value = '  +3.73%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.01%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '512.05'
$ws.Range("E30").Value = '  +2.98%  '
# Row 31
$ws.Range("D31").Value = '0.0₃0897'
$ws.Range("E31").Value = '  +0.03%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.81'
$ws.Range("E32").Value = '  +1.21%  '
# Row 33
$ws.Range("E33").Value = '  +2.52%  '
# Row 34
$ws.Range("E34").Value = '  +1.56%  '
# Row 35
$ws.Range("E35").Value = '  +0.03%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.09'
$ws.Range("E36").Value = '  +2.21%  '
# Row 37
$ws.Range("E37").Value = '  +0.39%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.42'
$ws.Range("E38").Value = '  +1.65%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.67'
$ws.Range("E39").Value = '  +0.01%  '
# Row 40
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.77'
$ws.Range("E40").Value = '  +5.60%  '
# Row 41
$ws.Range("B41").Value = 'ImmutableX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.32'
$ws.Range("E41").Value = '  -0.05%  '
# Row 42
$ws.Range("E42").Value = '  +0.04%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.85'
$ws.Range("E43").Value = '  +1.04%  '
# Row 44
$ws.Range("E44").Value = '  +0.25%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.42'
$ws.Range("E45").Value = '  +1.44%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '153.39'
$ws.Range("E46").Value = '  +7.73%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.57'
$ws.Range("E47").Value = '  +2.91%  '
# Row 48
$ws.Range("E48").Value = '  +2.89%  '
# Row 49
$ws.Range("D49").Value = '0.0₆0260'
$ws.Range("E49").Value = '  +2.89%  '
# Row 50
$ws.Range("E50").Value = '  +3.28%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0741'
$ws.Range("E51").Value = '  +0.18%  '
